# Edit script for "[Idea] In Eggcelent Condition.docx"
# Applies the changes described by the commit:
#   "Actual first map, main menu, outlines and decorations, egg cannon model"
#
# Strategy: operate on the document bottom-to-top (by paragraph index) so
# that earlier (lower-index) edits are not affected by index shifts caused
# by later (higher-index) deletions/insertions.

$d = $word.ActiveDocument

function Set-ParaXml {
    param($Paragraph, $InnerParaXml)
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $InnerParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Paragraph.Range.InsertXML($wrapper)
}

# ---------------------------------------------------------------------
# 1) "Egg Cannons" section near the end (original paragraphs 71-75)
#    - delete empty paragraph (71) + "Egg Cannons" heading (72)
#    - rewrite "Slowly aim in different directions." (73) -> new text, numId 5->21
#    - rewrite "Use a random force when shooting." (74) -> new text, numId 5->21
#    - delete "Create an actual 3D model..." (75)
# ---------------------------------------------------------------------

# Delete paragraph 75 ("Create an actual 3D model for it, plus a very fancy
# animation (windup and execute).")
$p75 = $d.Paragraphs(75)
$r = $d.Range($p75.Range.Start, $p75.Range.End)
$r.Delete()

# Rewrite paragraph 74 ("Use a random force when shooting.")
$p74 = $d.Paragraphs(74)
$xml74 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr>'
$xml74 = $xml74 + '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Egg cannons: </w:t></w:r>'
$xml74 = $xml74 + '<w:r><w:t>create fancy animation (windup + big shot)</w:t></w:r>'
$xml74 = $xml74 + '</w:p>'
Set-ParaXml $p74 $xml74

# Rewrite paragraph 73 ("Slowly aim in different directions.")
$p73 = $d.Paragraphs(73)
$xml73 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr>'
$xml73 = $xml73 + '<w:r><w:t xml:space="preserve">Try moving the second player with the </w:t></w:r>'
$xml73 = $xml73 + '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>right joystick</w:t></w:r>'
$xml73 = $xml73 + '<w:r><w:t xml:space="preserve"> on controller. (Would require adding an extra entry to the input map with some made-up number that will never be reached in real-life.)</w:t></w:r>'
$xml73 = $xml73 + '</w:p>'
Set-ParaXml $p73 $xml73

# Delete the empty paragraph (71) + "Egg Cannons" heading (72), keeping the
# (now-merged) empty paragraph in place.
$p71 = $d.Paragraphs(71)
$p72 = $d.Paragraphs(72)
$r = $d.Range($p71.Range.Start, $p72.Range.End)
$r.Delete()

# ---------------------------------------------------------------------
# 2) Remove stray <w:lastRenderedPageBreak/> on the "Show" run (paragraph 67)
# ---------------------------------------------------------------------
$p67 = $d.Paragraphs(67)
$xml67 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr>'
$xml67 = $xml67 + '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Show</w:t></w:r>'
$xml67 = $xml67 + '<w:r><w:t xml:space="preserve"> the area of influence of those powerups (attract, repel, freeze)</w:t></w:r>'
$xml67 = $xml67 + '</w:p>'
Set-ParaXml $p67 $xml67

# ---------------------------------------------------------------------
# 3) Delete the entire "Menu" section (heading 56 + bullets 57-61) plus the
#    trailing empty paragraph (62), keeping the empty paragraph (55) that
#    already precedes it (so exactly one blank line remains before "Fixes").
# ---------------------------------------------------------------------
$p56 = $d.Paragraphs(56)
$p62 = $d.Paragraphs(62)
$r = $d.Range($p56.Range.Start, $p62.Range.End)
$r.Delete()

# ---------------------------------------------------------------------
# 4) Rewrite the "Arenas" bullet (originally "Create training arena", 54)
# ---------------------------------------------------------------------
$p54 = $d.Paragraphs(54)
$xml54 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr>'
$xml54 = $xml54 + '<w:r><w:t xml:space="preserve">Training: </w:t></w:r>'
$xml54 = $xml54 + '<w:r><w:t>Two cannons, with some light random swiveling</w:t></w:r>'
$xml54 = $xml54 + '<w:r><w:br/><w:t xml:space="preserve">Training: </w:t></w:r>'
$xml54 = $xml54 + '<w:r><w:t xml:space="preserve">some egg-shaped rocks or statues? </w:t></w:r>'
$xml54 = $xml54 + '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>A splash of color?</w:t></w:r>'
$xml54 = $xml54 + '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> =&gt; When an egg breaks on the ground, it paints something into the texture used for the outline overlay!</w:t></w:r>'
$xml54 = $xml54 + '</w:p>'
Set-ParaXml $p54 $xml54

# ---------------------------------------------------------------------
# 5) Insert a new empty paragraph right after "Implement the last ones" (52)
#    and before the "Arenas" heading (53).
# ---------------------------------------------------------------------
$p52 = $d.Paragraphs(52)
$insertionPoint = $d.Range($p52.Range.End, $p52.Range.End)
$xmlNewPara = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xmlNewPara)

# ---------------------------------------------------------------------
# 6) Delete the entire "Essentials" section (heading 43 + bullets 44-49 +
#    trailing empty paragraph 50), leaving the "Eggs" heading (51) as the
#    first Heading2 under "To Do".
# ---------------------------------------------------------------------
$p43 = $d.Paragraphs(43)
$p50 = $d.Paragraphs(50)
$r = $d.Range($p43.Range.Start, $p50.Range.End)
$r.Delete()
